$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 29.333334
$ws.Range("I8").Value = 26.3
$ws.Range("J8").Value = 44.5
$ws.Range("K8").Value = 78.90000000000001
$ws.Range("L8").Value = 133.5
$ws.Range("M8").Value = 60.09999999999999
$ws.Range("N8").Value = -411.5

$ws.Range("H18").Value = 305.05554
$ws.Range("I18").Value = 305.05554
$ws.Range("K18").Value = 305.05554
$ws.Range("M18").Value = -21.05554000000001

$ws.Range("H19").Value = 1406.7
$ws.Range("I19").Value = 1899.6154
$ws.Range("J19").Value = 491.2857
$ws.Range("K19").Value = 1899.6154
$ws.Range("L19").Value = 491.2857
$ws.Range("M19").Value = -1724.6154
$ws.Range("N19").Value = -841.2857

$ws.Range("H33").Value = 9210527
$ws.Range("I33").Value = 12558119
$ws.Range("J33").Value = 4647.75
$ws.Range("K33").Value = 12558119
$ws.Range("L33").Value = 4647.75
$ws.Range("M33").Value = -12557890
$ws.Range("N33").Value = -5105.75

$ws.Range("H40").Value = 1941.7693
$ws.Range("I40").Value = 1799.8
$ws.Range("J40").Value = 2030.5
$ws.Range("K40").Value = 1799.8
$ws.Range("L40").Value = 2030.5
$ws.Range("M40").Value = -1624.8
$ws.Range("N40").Value = -2380.5

$ws.Range("H58").Value = 1496
$ws.Range("I58").Value = 994.6667
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 2984.0001
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -2834.0001
$ws.Range("N58").Value = -9300

$ws.Range("H80").Value = 668076.5600000001
$ws.Range("I80").Value = 1171.3334
$ws.Range("J80").Value = 1112680
$ws.Range("K80").Value = 3514.0002
$ws.Range("L80").Value = 3338040
$ws.Range("M80").Value = -2516.0002
$ws.Range("N80").Value = -3340036

$ws.Range("H83").Value = 668076.5600000001
$ws.Range("I83").Value = 1171.3334
$ws.Range("J83").Value = 1112680
$ws.Range("K83").Value = 10542.0006
$ws.Range("L83").Value = 10014120
$ws.Range("M83").Value = -5550.000599999999
$ws.Range("N83").Value = -10024104

$ws.Range("H113").Value = 4781.923
$ws.Range("I113").Value = 3520.375
$ws.Range("J113").Value = 6800.4
$ws.Range("K113").Value = 3520.375
$ws.Range("L113").Value = 6800.4
$ws.Range("M113").Value = -266.375
$ws.Range("N113").Value = -13308.4

$ws.Range("H116").Value = 7128.4614
$ws.Range("J116").Value = 9647.666999999999
$ws.Range("L116").Value = 9647.666999999999
$ws.Range("N116").Value = -16531.667

$ws.Range("H132").Value = 1237.9231
$ws.Range("I132").Value = 1248.7142
$ws.Range("J132").Value = 1143.5
$ws.Range("K132").Value = 3746.1426
$ws.Range("L132").Value = 3430.5
$ws.Range("M132").Value = -1216.1426
$ws.Range("N132").Value = -8490.5

$ws.Range("H135").Value = 2843.4167
$ws.Range("I135").Value = 2374.6365
$ws.Range("K135").Value = 21371.7285
$ws.Range("M135").Value = -18836.7285

$ws.Range("H137").Value = 11600.2
$ws.Range("I137").Value = 4941.1665
$ws.Range("K137").Value = 14823.4995
$ws.Range("M137").Value = -12273.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4377.7915
$ws.Range("I2").Value = 3741.3333
$ws.Range("J2").Value = 5438.5557
$ws.Range("K2").Value = 3741.3333
$ws.Range("L2").Value = 5438.5557
$ws.Range("M2").Value = -3628.3333
$ws.Range("N2").Value = -5664.5557

$ws.Range("H32").Value = 3107.6584
$ws.Range("I32").Value = 3107.6584
$ws.Range("K32").Value = 3107.6584
$ws.Range("M32").Value = -2820.6584

$ws.Range("H74").Value = 20277.238
$ws.Range("I74").Value = 21462.334
$ws.Range("K74").Value = 21462.334
$ws.Range("M74").Value = -20588.334

$ws.Range("H77").Value = 20277.238
$ws.Range("I77").Value = 21462.334
$ws.Range("K77").Value = 107311.67
$ws.Range("M77").Value = -102943.67

$ws.Range("H116").Value = 4377.7915
$ws.Range("I116").Value = 3741.3333
$ws.Range("J116").Value = 5438.5557
$ws.Range("K116").Value = 3741.3333
$ws.Range("L116").Value = 5438.5557
$ws.Range("M116").Value = -1447.3333
$ws.Range("N116").Value = -10026.5557

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4377.7915
$ws.Range("I3").Value = 3741.3333
$ws.Range("J3").Value = 5438.5557
$ws.Range("K3").Value = 3741.3333
$ws.Range("L3").Value = 5438.5557
$ws.Range("M3").Value = -3627.3333
$ws.Range("N3").Value = -5666.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2421.5881
$ws.Range("I16").Value = 2497.6428
$ws.Range("K16").Value = 2497.6428
$ws.Range("M16").Value = -2210.6428

$ws.Range("H22").Value = 430.8
$ws.Range("I22").Value = 198.5
$ws.Range("J22").Value = 1360
$ws.Range("K22").Value = 198.5
$ws.Range("L22").Value = 1360
$ws.Range("M22").Value = 151.5
$ws.Range("N22").Value = -2060

$ws.Range("H23").Value = 26293.5
$ws.Range("J23").Value = 27218
$ws.Range("L23").Value = 27218
$ws.Range("N23").Value = -27698

$ws.Range("H27").Value = 26293.5
$ws.Range("J27").Value = 27218
$ws.Range("L27").Value = 27218
$ws.Range("N27").Value = -27602

$ws.Range("H58").Value = 4326.2354
$ws.Range("I58").Value = 2831.9092
$ws.Range("J58").Value = 5040.913
$ws.Range("K58").Value = 2831.9092
$ws.Range("L58").Value = 5040.913
$ws.Range("M58").Value = -2628.9092
$ws.Range("N58").Value = -5446.913

$ws.Range("H105").Value = 2426.0667
$ws.Range("I105").Value = 2337.8462
$ws.Range("J105").Value = 2999.5
$ws.Range("K105").Value = 2337.8462
$ws.Range("L105").Value = 2999.5
$ws.Range("M105").Value = -590.8462
$ws.Range("N105").Value = -6493.5

$ws.Range("H113").Value = 2421.5881
$ws.Range("I113").Value = 2497.6428
$ws.Range("K113").Value = 2497.6428
$ws.Range("M113").Value = -327.6428000000001

$ws.Range("H134").Value = 3914.2632
$ws.Range("I134").Value = 2295.4546
$ws.Range("K134").Value = 6886.3638
$ws.Range("M134").Value = -4351.3638

$ws.Range("H136").Value = 4326.2354
$ws.Range("I136").Value = 2831.9092
$ws.Range("J136").Value = 5040.913
$ws.Range("K136").Value = 8495.7276
$ws.Range("L136").Value = 15122.739
$ws.Range("M136").Value = -5945.7276
$ws.Range("N136").Value = -20222.739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 56.64706
$ws.Range("J2").Value = 79.59999999999999
$ws.Range("L2").Value = 477.6
$ws.Range("N2").Value = -703.5999999999999

$ws.Range("H4").Value = 49189990
$ws.Range("I4").Value = 49189990
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 147569970
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -147569858
$ws.Range("N4").ClearContents()

$ws.Range("H7").Value = 71508300
$ws.Range("I7").Value = 166667000
$ws.Range("K7").Value = 500001000
$ws.Range("M7").Value = -500000888

$ws.Range("H15").Value = 27.666666
$ws.Range("I15").Value = 33.166668
$ws.Range("K15").Value = 99.500004
$ws.Range("M15").Value = 40.499996

$ws.Range("H20").Value = 1998
$ws.Range("J20").Value = 1998
$ws.Range("L20").Value = 5994
$ws.Range("N20").Value = -6448

$ws.Range("H22").Value = 67300.60000000001
$ws.Range("I22").Value = 390
$ws.Range("J22").Value = 167666.5
$ws.Range("K22").Value = 1170
$ws.Range("L22").Value = 502999.5
$ws.Range("M22").Value = -1001
$ws.Range("N22").Value = -503337.5

$ws.Range("H27").Value = 67300.60000000001
$ws.Range("I27").Value = 390
$ws.Range("J27").Value = 167666.5
$ws.Range("K27").Value = 1170
$ws.Range("L27").Value = 502999.5
$ws.Range("M27").Value = -1068
$ws.Range("N27").Value = -503203.5

$ws.Range("H44").Value = 51996.5
$ws.Range("J44").Value = 3994
$ws.Range("L44").Value = 11982
$ws.Range("N44").Value = -12778

$ws.Range("H75").Value = 66667456
$ws.Range("J75").Value = 66667456
$ws.Range("L75").Value = 200002368
$ws.Range("N75").Value = -200004364

$ws.Range("H78").Value = 66667456
$ws.Range("J78").Value = 66667456
$ws.Range("L78").Value = 600007104
$ws.Range("N78").Value = -600017088

$ws.Range("H107").Value = 884.6
$ws.Range("I107").Value = 690
$ws.Range("J107").Value = 1176.5
$ws.Range("K107").Value = 2070
$ws.Range("L107").Value = 3529.5
$ws.Range("M107").Value = -150
$ws.Range("N107").Value = -7369.5

$ws.Range("H113").Value = 1112.7
$ws.Range("J113").Value = 840.8333
$ws.Range("L113").Value = 2522.4999
$ws.Range("N113").Value = -6862.4999

$ws.Range("H122").Value = 7143543
$ws.Range("J122").Value = 11111947
$ws.Range("L122").Value = 100007523
$ws.Range("N122").Value = -100012423

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4859.1875
$ws.Range("I100").Value = 3439.1428
$ws.Range("K100").Value = 3439.1428
$ws.Range("M100").Value = -2898.1428

$ws.Range("H122").Value = 3280.913
$ws.Range("I122").Value = 3823.4443
$ws.Range("J122").Value = 2932.1428
$ws.Range("K122").Value = 11470.3329
$ws.Range("L122").Value = 8796.428400000001
$ws.Range("M122").Value = -9020.332900000001
$ws.Range("N122").Value = -13696.4284

$ws.Range("H136").Value = 5360
$ws.Range("J136").Value = 5954.875
$ws.Range("L136").Value = 17864.625
$ws.Range("N136").Value = -22964.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5557059.5
$ws.Range("I107").Value = 960.2
$ws.Range("J107").Value = 22225358
$ws.Range("K107").Value = 2880.6
$ws.Range("L107").Value = 66676074
$ws.Range("M107").Value = -960.6000000000004
$ws.Range("N107").Value = -66679914

$ws.Range("H113").Value = 661.1429000000001
$ws.Range("J113").Value = 1470
$ws.Range("L113").Value = 4410
$ws.Range("N113").Value = -8750

$ws.Range("H132").Value = 141166.1
$ws.Range("I132").Value = 260260.84
$ws.Range("K132").Value = 780782.52
$ws.Range("M132").Value = -778252.52

$ws.Range("H136").Value = 8336484.5
$ws.Range("I136").Value = 25002888
$ws.Range("J136").Value = 3283.4375
$ws.Range("K136").Value = 75008664
$ws.Range("L136").Value = 9850.3125
$ws.Range("M136").Value = -75006114
$ws.Range("N136").Value = -14950.3125
